$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '22.080.63'
$ws.Cells.Item(2, 5).Value = '  -1.51%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.556.57'
$ws.Cells.Item(3, 5).Value = '  -0.74%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9982'
$ws.Cells.Item(4, 5).Value = '  -0.24%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.9989'
$ws.Cells.Item(5, 5).Value = '  -0.16%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '287.19'
$ws.Cells.Item(6, 5).Value = '  +0.06%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.3846'
$ws.Cells.Item(7, 5).Value = '  +3.53%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3250'
$ws.Cells.Item(8, 5).Value = '  -1.51%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '41.43'
$ws.Cells.Item(9, 5).Value = '  -12.09%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.126'
$ws.Cells.Item(10, 5).Value = '  -2.66%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07331'
$ws.Cells.Item(11, 5).Value = '  -1.63%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.9982'
$ws.Cells.Item(12, 5).Value = '  -0.27%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '19.47'
$ws.Cells.Item(13, 5).Value = '  -5.62%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.726'
$ws.Cells.Item(14, 5).Value = '  -2.57%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.812'
$ws.Cells.Item(15, 5).Value = '  -0.52%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '1.555.72'
$ws.Cells.Item(16, 5).Value = '  +0.00%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.00001092'
$ws.Cells.Item(17, 5).Value = '  -1.45%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.06622'
$ws.Cells.Item(18, 5).Value = '  -1.12%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '85.17'
$ws.Cells.Item(19, 5).Value = '  -1.77%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.408'
$ws.Cells.Item(20, 5).Value = '  +0.72%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.9992'
$ws.Cells.Item(21, 5).Value = '  -0.07%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '15.98'
$ws.Cells.Item(22, 5).Value = '  -2.59%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '11.48'
$ws.Cells.Item(23, 5).Value = '  -3.42%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '22.091.60'
$ws.Cells.Item(24, 5).Value = '  -1.38%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -1.13%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.547'
$ws.Cells.Item(26, 5).Value = '  -1.69%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -1.61%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '18.90'
$ws.Cells.Item(28, 5).Value = '  -2.80%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '4.853'
$ws.Cells.Item(29, 5).Value = '  -1.77%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '1.728.42'
$ws.Cells.Item(30, 5).Value = '  -0.41%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '120.81'
$ws.Cells.Item(31, 5).Value = '  -2.72%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.103'
$ws.Cells.Item(32, 5).Value = '  +4.00%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.903'
$ws.Cells.Item(33, 5).Value = '  -2.05%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.667'
$ws.Cells.Item(34, 5).Value = '  -15.48%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'FraxShare'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '9.277'
$ws.Cells.Item(35, 5).Value = '  -5.17%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Stellar'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.08152'
$ws.Cells.Item(36, 5).Value = '  -1.51%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'VeChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.02300'
$ws.Cells.Item(37, 5).Value = '  -4.75%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '5.255'
$ws.Cells.Item(38, 5).Value = '  -0.42%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.06200'
$ws.Cells.Item(39, 5).Value = '  -2.57%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.2110'
$ws.Cells.Item(40, 5).Value = '  -4.00%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.221'
$ws.Cells.Item(41, 5).Value = '  -5.19%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '10.90'
$ws.Cells.Item(42, 5).Value = '  -3.64%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.9990'
$ws.Cells.Item(43, 5).Value = '  -0.07%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.5948'
$ws.Cells.Item(44, 5).Value = '  -3.40%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.50'
$ws.Cells.Item(45, 5).Value = '  -2.01%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.720'
$ws.Cells.Item(46, 5).Value = '  -1.00%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.5761'
$ws.Cells.Item(47, 5).Value = '  -3.72%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.937'
$ws.Cells.Item(48, 5).Value = '  -4.66%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '119.69'
$ws.Cells.Item(49, 5).Value = '  -3.44%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.156'
$ws.Cells.Item(50, 5).Value = '  -3.02%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06895'
$ws.Cells.Item(51, 5).Value = '  -3.89%  '
